$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.52000000000039
$ws.Range("H2").Value = [double]"2.259995978880726e-16"
$ws.Range("K2").Value = 30.7447986016462
$ws.Range("L2").Value = "[23.718961859237965, 37.77063534405444]"
$ws.Range("M2").Value = [double]"2.220446049250313e-16"
$ws.Range("N2").Value = [double]"2.220446049250313e-16"
$ws.Range("O2").Value = 1.817658211986887
$ws.Range("P2").Value = "[1.553500271144502, 2.0818161528292727]"
$ws.Range("S2").Value = 54.76603129296286
$ws.Range("T2").Value = "[50.628822282845135, 58.90324030308058]"
$ws.Range("W2").Value = 17.42662662662691
$ws.Range("X2").Value = 16.39575575575602
$ws.Range("Y2").Value = 18.45749749749779

# Row 3 updates
$ws.Range("E3").Value = 25.09000000000048
$ws.Range("H3").Value = [double]"2.259995978880726e-16"
$ws.Range("K3").Value = 44.51316611264501
$ws.Range("L3").Value = "[36.03426077116279, 52.99207145412723]"
$ws.Range("O3").Value = -2.943474197958005
$ws.Range("P3").Value = "[-3.157316340544698, -2.7296320553713116]"
$ws.Range("S3").Value = 59.74528266148919
$ws.Range("T3").Value = "[54.61722966932098, 64.8733356536574]"
$ws.Range("W3").Value = 11.7538738738741
$ws.Range("X3").Value = 10.89995995996017
$ws.Range("Y3").Value = 12.60778778778804
